$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point rounding in A8 (re-saved timestamp value)
$ws.Cells.Item(8, 1).Value = 45864.50030652778

# Add new row 9 with the latest sensor reading pulled in by the scheduled task
$ws.Cells.Item(9, 1).Value = 45864.54198652151
$ws.Cells.Item(9, 1).NumberFormat = $ws.Cells.Item(8, 1).NumberFormat

$ws.Cells.Item(9, 2).Value = 2025
$ws.Cells.Item(9, 3).Value = 30
$ws.Cells.Item(9, 4).Value = 21.11
$ws.Cells.Item(9, 5).Value = 68.28
$ws.Cells.Item(9, 6).Value = 611.72
$ws.Cells.Item(9, 7).Value = 12.73
$ws.Cells.Item(9, 8).Value = "ESE"
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = "13:00:27"
